$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.762.75'
$ws.Range('E2').Value = '  -2.41%  '
$ws.Range('D3').Value = '3.140.78'
$ws.Range('E3').Value = '  -2.38%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.05'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.06'
$ws.Range('E6').Value = '  -5.45%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.129.79'
$ws.Range('E8').Value = '  -2.70%  '
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.150'
$ws.Range('E10').Value = '  -4.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.29'
$ws.Range('E11').Value = '  -4.66%  '
$ws.Range('E12').Value = '  -2.16%  '
$ws.Range('E13').Value = '  -2.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.44'
$ws.Range('E14').Value = '  -6.00%  '
$ws.Range('D15').Value = '3.652.79'
$ws.Range('E15').Value = '  -3.06%  '
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = '63.809.07'
$ws.Range('E17').Value = '  -2.81%  '
$ws.Range('D18').Value = '3.134.26'
$ws.Range('E18').Value = '  -3.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.83'
$ws.Range('E19').Value = '  -3.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '476.83'
$ws.Range('E20').Value = '  -2.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.61'
$ws.Range('E21').Value = '  -2.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.714'
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.82'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.58'
$ws.Range('E24').Value = '  -4.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.64'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('E27').Value = '  -5.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.49'
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.123'
$ws.Range('E29').Value = '  -5.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  +1.77%  '
$ws.Range('E31').Value = '  -10.13%  '
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.70'
$ws.Range('E33').Value = '  -2.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.22'
$ws.Range('E34').Value = '  -3.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.12'
$ws.Range('E35').Value = '  +0.98%  '
$ws.Range('D36').Value = '0.0₃0792'
$ws.Range('E36').Value = '  +7.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.94'
$ws.Range('E37').Value = '  -4.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '52.84'
$ws.Range('E38').Value = '  -4.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '458.78'
$ws.Range('E39').Value = '  -4.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.03'
$ws.Range('E40').Value = '  -8.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0394'
$ws.Range('E41').Value = '  -4.27%  '
$ws.Range('E42').Value = '  -6.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.33'
$ws.Range('E43').Value = '  -2.76%  '
$ws.Range('D44').Value = '2.848.89'
$ws.Range('E44').Value = '  -3.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.30'
$ws.Range('E45').Value = '  -6.89%  '
$ws.Range('E46').Value = '  -5.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.43'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.30'
$ws.Range('E48').Value = '  -5.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '118.99'
$ws.Range('E51').Value = '  -1.84%  '
